$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the typo "Mortize locks" -> "Mortise locks" in C5
$ws.Range("C5").Value = "Mortise locks"

# Best-fit / autofit column C to its (now updated) contents
$ws.Columns("C:C").AutoFit()

# Update the active selection from E10 to F9
$ws.Range("F9").Select()
